$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were deleted from the source data (RM 232, SC 92).
# Deleting row 26 first, then row 27 (which is "SC 92" after the first deletion
# shifts everything up by one), shifts all subsequent rows up by two in total,
# matching the new dimension A1:F33.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply the individual cell value corrections.
$ws.Range("F2").Value = $null
$ws.Range("C3").Value = 11.2
$ws.Range("F3").Value = $null
$ws.Range("D4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = 17.66
$ws.Range("F8").Value = $null
$ws.Range("D9").Value = -14.5
$ws.Range("F9").Value = 17.26
$ws.Range("D10").Value = -14.7
$ws.Range("D11").Value = -15.5
$ws.Range("D12").Value = -14.1
$ws.Range("E12").Value = $null
$ws.Range("E13").Value = -5.3
$ws.Range("F13").Value = $null
$ws.Range("E14").Value = -5.4
$ws.Range("D15").Value = $null
$ws.Range("F15").Value = 16.2
$ws.Range("D17").Value = $null
$ws.Range("E17").Value = -7.3
$ws.Range("F17").Value = $null
$ws.Range("D18").Value = $null
$ws.Range("F19").Value = 17.81
$ws.Range("D20").Value = $null
$ws.Range("E20").Value = $null
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = $null
$ws.Range("E23").Value = $null
$ws.Range("F23").Value = $null
$ws.Range("E25").Value = $null
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = $null
$ws.Range("E28").Value = -5.9
$ws.Range("D31").Value = -13.7
$ws.Range("F31").Value = 17.18
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
